# Edit script: "Add title page and table of contents"
# Applies the geometry / text / shape changes that correspond to the
# single slide contained in before.pptx (the flow-chart slide, whose
# p14:creationId is 3530386684 in the original multi-slide deck).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1. "Rectangle 2" (id=4): resize/reposition and give it the caption
#    "Argument parser".
# ---------------------------------------------------------------------
$sh4 = $s.Shapes.Item("Rectangle 2")
$sh4.Left   = 1835696 / 12700
$sh4.Top    = 2564904 / 12700
$sh4.Width  = 1337272 / 12700
$sh4.Height = 288032 / 12700
$sh4.TextFrame.TextRange.Text = "Argument parser"
$sh4.TextFrame.TextRange.Font.Size = 12

# ---------------------------------------------------------------------
# 2. Small positional nudges (x-shift mostly) of several existing
#    shapes, to make room for the redrawn flow-chart.
# ---------------------------------------------------------------------
$sh5 = $s.Shapes.Item("AutoShape 3")
$sh5.Left = 5172695 / 12700

$sh6 = $s.Shapes.Item("AutoShape 4")
$sh6.Left = 2737587 / 12700

$sh7 = $s.Shapes.Item("AutoShape 5")
$sh7.Left = 3787421 / 12700

$sh1031 = $s.Shapes.Item("AutoShape 7")
$sh1031.Left = 3303531 / 12700

$sh13 = $s.Shapes.Item("Правая фигурная скобка 12")
$sh13.Left = 971600 / 12700
$sh13.Top  = 2060848 / 12700

# ---------------------------------------------------------------------
# 3. Delete "TextBox 14" ("Входные данные") - replaced further down by
#    the new "Input data" oval (id 58).
# ---------------------------------------------------------------------
$s.Shapes.Item("TextBox 14").Delete()

# ---------------------------------------------------------------------
# 4. "Прямая со стрелкой 16" (id=17): reposition + flip vertically.
# ---------------------------------------------------------------------
$sh17 = $s.Shapes.Item("Прямая со стрелкой 16")
$sh17.Left   = 1403648 / 12700
$sh17.Top    = 2708920 / 12700
$sh17.Width  = 432048 / 12700
$sh17.Height = 17381 / 12700
$sh17.Flip(1)

# ---------------------------------------------------------------------
# 5. "TextBox 19" (id=20): reposition.
# ---------------------------------------------------------------------
$sh20 = $s.Shapes.Item("TextBox 19")
$sh20.Left = 3347864 / 12700
$sh20.Top  = 2492896 / 12700

# ---------------------------------------------------------------------
# 6. Further x-shift nudges.
# ---------------------------------------------------------------------
$sh31 = $s.Shapes.Item("Скругленный прямоугольник 30")
$sh31.Left = 3945131 / 12700

$sh34 = $s.Shapes.Item("Соединительная линия уступом 33")
$sh34.Left = 5083565 / 12700

$sh38 = $s.Shapes.Item("Соединительная линия уступом 37")
$sh38.Left = 3869475 / 12700

$sh40 = $s.Shapes.Item("Соединительная линия уступом 39")
$sh40.Left = 4788025 / 12700

$sh53 = $s.Shapes.Item("Соединительная линия уступом 52")
$sh53.Left   = 3172968 / 12700
$sh53.Top    = 2708920 / 12700
$sh53.Width  = 1262525 / 12700
$sh53.Height = 220870 / 12700

# ---------------------------------------------------------------------
# 7. New shapes added at the end of the flow-chart.
# ---------------------------------------------------------------------

# 7a. "Прямоугольник 15" (id=16) - "Format the results and writing to file"
$sh16 = $s.Shapes.AddShape(1, 1691680/12700, 4797152/12700, 1296144/12700, 648072/12700)
$sh16.Name = "Прямоугольник 15"
$sh16.Line.Weight = 0.25
$sh16.Line.ForeColor.RGB = 0
$sh16.TextFrame.TextRange.Text = "Format the results and writing to file"
$sh16.TextFrame.TextRange.Font.Size = 12
$sh16.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# 7b. "Прямая со стрелкой 18" (id=19) connecting 31 -> 16
$sh19 = $s.Shapes.AddConnector(1, 0, 0, 1, 1)
$sh19.Name = "Прямая со стрелкой 18"
$sh19.Left   = 2987824 / 12700
$sh19.Top    = 5110014 / 12700
$sh19.Width  = 957307 / 12700
$sh19.Height = 11174 / 12700
$sh19.Flip(0)
$sh19.Line.ForeColor.RGB = 0
$sh19.Line.EndArrowheadStyle = 2

# 7c. "Овал 42" (id=43) - "File"
$sh43 = $s.Shapes.AddShape(9, 179512/12700, 4869160/12700, 864096/12700, 576064/12700)
$sh43.Name = "Овал 42"
$sh43.Line.Weight = 0.25
$sh43.Line.ForeColor.RGB = 0
$sh43.TextFrame.TextRange.Text = "File"
$sh43.TextFrame.TextRange.Font.Size = 12
$sh43.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# 7d. "Прямая со стрелкой 44" (id=45) connecting 16 -> 43
$sh45 = $s.Shapes.AddConnector(1, 0, 0, 1, 1)
$sh45.Name = "Прямая со стрелкой 44"
$sh45.Left   = 1043608 / 12700
$sh45.Top    = 5121188 / 12700
$sh45.Width  = 648072 / 12700
$sh45.Height = 36004 / 12700
$sh45.Flip(0)
$sh45.Line.ForeColor.RGB = 0
$sh45.Line.EndArrowheadStyle = 2

# 7e. "Прямая соединительная линия 46" (id=47) - dashed vertical divider
$sh47 = $s.Shapes.AddLine(1475656/12700, 2348880/12700, 1475656/12700, (2348880+3240360)/12700)
$sh47.Name = "Прямая соединительная линия 46"
$sh47.Line.Weight = 1
$sh47.Line.DashStyle = 4
$sh47.Line.ForeColor.RGB = 0

# 7f. "Овал 57" (id=58) - "Input data"
$sh58 = $s.Shapes.AddShape(9, 107504/12700, 2420888/12700, 1008112/12700, 648072/12700)
$sh58.Name = "Овал 57"
$sh58.Line.Weight = 0.25
$sh58.Line.ForeColor.RGB = 0
$sh58.TextFrame.TextRange.Text = "`rInput data`r"
$sh58.TextFrame.TextRange.Font.Size = 12
$sh58.TextFrame.TextRange.Paragraphs(2).ParagraphFormat.Alignment = 2
$sh58.TextFrame.TextRange.Paragraphs(3).ParagraphFormat.Alignment = 2
